$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, [string]$Text)
    $looksNumeric = $Text -match '^\s*[+-]?\d+(\.\d+)?\s*$'
    if ($looksNumeric) {
        $Cell.NumberFormat = "@"
        $Cell.Value = $Text
        $Cell.Style = "Normal"
    } else {
        $Cell.Value = $Text
    }
}

# Row 2
Set-TextValue $ws.Range("D2") '25.975.85'
Set-TextValue $ws.Range("E2") '  -0.02%  '

# Row 3
Set-TextValue $ws.Range("D3") '1.633.54'
Set-TextValue $ws.Range("E3") '  -0.48%  '

# Row 4
Set-TextValue $ws.Range("E4") '  +0.10%  '

# Row 5
Set-TextValue $ws.Range("D5") '213.90'
Set-TextValue $ws.Range("E5") '  -1.02%  '

# Row 6
Set-TextValue $ws.Range("E6") '  -0.64%  '

# Row 7
Set-TextValue $ws.Range("E7") '  +0.15%  '

# Row 8
Set-TextValue $ws.Range("E8") '  -1.66%  '

# Row 9
Set-TextValue $ws.Range("E9") '  -2.43%  '

# Row 10
Set-TextValue $ws.Range("E10") '  -5.72%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.0792'
Set-TextValue $ws.Range("E11") '  -0.33%  '

# Row 12
Set-TextValue $ws.Range("D12") '1.859.96'

# Row 13
Set-TextValue $ws.Range("B13") 'Polkadot'
Set-TextValue $ws.Range("C13") 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range("D13") '4.20'
Set-TextValue $ws.Range("E13") '  -2.13%  '

# Row 14
Set-TextValue $ws.Range("B14") 'WrappedEther'
Set-TextValue $ws.Range("C14") 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D14") '1.632.42'
Set-TextValue $ws.Range("E14") '  -0.24%  '

# Row 16
Set-TextValue $ws.Range("E16") '  -2.66%  '

# Row 17
Set-TextValue $ws.Range("D17") '25.979.82'
Set-TextValue $ws.Range("E17") '  +0.09%  '

# Row 18
Set-TextValue $ws.Range("D18") '61.74'

# Row 19
Set-TextValue $ws.Range("E19") '  +0.12%  '

# Row 20
Set-TextValue $ws.Range("D20") '190.39'
Set-TextValue $ws.Range("E20") '  -1.44%  '

# Row 21
Set-TextValue $ws.Range("E21") '  -2.73%  '

# Row 22
Set-TextValue $ws.Range("E22") '  -3.65%  '

# Row 23
Set-TextValue $ws.Range("D23") '6.13'
Set-TextValue $ws.Range("E23") '  -1.89%  '

# Row 24
Set-TextValue $ws.Range("D24") '0.133'
Set-TextValue $ws.Range("E24") '  +0.08%  '

# Row 25
Set-TextValue $ws.Range("D25") '143.48'
Set-TextValue $ws.Range("E25") '  -0.67%  '

# Row 26
Set-TextValue $ws.Range("E26") '  +0.11%  '

# Row 27
Set-TextValue $ws.Range("E27") '  -3.13%  '

# Row 28
Set-TextValue $ws.Range("E28") '  -1.99%  '

# Row 29
Set-TextValue $ws.Range("D29") '15.22'
Set-TextValue $ws.Range("E29") '  -2.03%  '

# Row 30
Set-TextValue $ws.Range("E30") '  -1.35%  '

# Row 31
Set-TextValue $ws.Range("E31") '  -3.05%  '

# Row 32
Set-TextValue $ws.Range("D32") '3.15'
Set-TextValue $ws.Range("E32") '  -2.98%  '

# Row 33
Set-TextValue $ws.Range("E33") '  -4.16%  '

# Row 34
Set-TextValue $ws.Range("E34") '  -1.89%  '

# Row 35
Set-TextValue $ws.Range("E35") '  -2.60%  '

# Row 36
Set-TextValue $ws.Range("D36") '1.134.63'
Set-TextValue $ws.Range("E36") '  +0.14%  '

# Row 37
Set-TextValue $ws.Range("D37") '0.865'
Set-TextValue $ws.Range("E37") '  -4.40%  '

# Row 38
Set-TextValue $ws.Range("E38") '  -1.26%  '

# Row 39
Set-TextValue $ws.Range("E39") '  -3.09%  '

# Row 40
Set-TextValue $ws.Range("E40") '  -0.98%  '

# Row 41
Set-TextValue $ws.Range("D41") '98.62'
Set-TextValue $ws.Range("E41") '  -0.63%  '

# Row 42
Set-TextValue $ws.Range("E42") '  -2.24%  '

# Row 43
Set-TextValue $ws.Range("E43") '  -5.06%  '

# Row 44
Set-TextValue $ws.Range("D44") '1.769.62'
Set-TextValue $ws.Range("E44") '  -0.44%  '

# Row 45
Set-TextValue $ws.Range("E45") '  -1.67%  '

# Row 46
Set-TextValue $ws.Range("D46") '55.11'
Set-TextValue $ws.Range("E46") '  -2.58%  '

# Row 47
Set-TextValue $ws.Range("E47") '  -0.72%  '

# Row 48
Set-TextValue $ws.Range("D48") '1.48'
Set-TextValue $ws.Range("E48") '  +1.47%  '

# Row 49
Set-TextValue $ws.Range("D49") '0.415'
Set-TextValue $ws.Range("E49") '  -0.15%  '

# Row 50
Set-TextValue $ws.Range("D50") '7.52'
Set-TextValue $ws.Range("E50") '  -2.31%  '

# Row 51
Set-TextValue $ws.Range("E51") '  +0.09%  '
